$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-07-22 Saturday", $true, $false, $false, $false, $false, $true, 0, $false, "2023-07-23 Sunday", 1) | Out-Null

# Update the multiplication table cells (row-major order, 20 rows x 5 cols)
# Using wdFindStop (0) for Wrap and wdReplaceOne (1) for Replace so each call
# only touches a single occurrence, applied in document order -- this matters
# because one pair of cells shares identical old text ("26x74=").
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("82×36=", $true, $false, $false, $false, $false, $true, 0, $false, "97×46=", 1) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("52×26=", $true, $false, $false, $false, $false, $true, 0, $false, "15×64=", 1) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("38×67=", $true, $false, $false, $false, $false, $true, 0, $false, "10×94=", 1) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("22×99=", $true, $false, $false, $false, $false, $true, 0, $false, "32×49=", 1) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("86×92=", $true, $false, $false, $false, $false, $true, 0, $false, "38×75=", 1) | Out-Null
$t.Cell(2, 1).Range.Find.Execute("57×89=", $true, $false, $false, $false, $false, $true, 0, $false, "43×89=", 1) | Out-Null
$t.Cell(2, 2).Range.Find.Execute("43×17=", $true, $false, $false, $false, $false, $true, 0, $false, "70×14=", 1) | Out-Null
$t.Cell(2, 3).Range.Find.Execute("26×74=", $true, $false, $false, $false, $false, $true, 0, $false, "38×47=", 1) | Out-Null
$t.Cell(2, 4).Range.Find.Execute("96×26=", $true, $false, $false, $false, $false, $true, 0, $false, "38×46=", 1) | Out-Null
$t.Cell(2, 5).Range.Find.Execute("39×47=", $true, $false, $false, $false, $false, $true, 0, $false, "33×74=", 1) | Out-Null
$t.Cell(3, 1).Range.Find.Execute("38×45=", $true, $false, $false, $false, $false, $true, 0, $false, "69×43=", 1) | Out-Null
$t.Cell(3, 2).Range.Find.Execute("30×46=", $true, $false, $false, $false, $false, $true, 0, $false, "98×18=", 1) | Out-Null
$t.Cell(3, 3).Range.Find.Execute("80×88=", $true, $false, $false, $false, $false, $true, 0, $false, "98×32=", 1) | Out-Null
$t.Cell(3, 4).Range.Find.Execute("26×84=", $true, $false, $false, $false, $false, $true, 0, $false, "42×18=", 1) | Out-Null
$t.Cell(3, 5).Range.Find.Execute("57×73=", $true, $false, $false, $false, $false, $true, 0, $false, "23×18=", 1) | Out-Null
$t.Cell(4, 1).Range.Find.Execute("53×22=", $true, $false, $false, $false, $false, $true, 0, $false, "32×54=", 1) | Out-Null
$t.Cell(4, 2).Range.Find.Execute("87×43=", $true, $false, $false, $false, $false, $true, 0, $false, "80×16=", 1) | Out-Null
$t.Cell(4, 3).Range.Find.Execute("81×32=", $true, $false, $false, $false, $false, $true, 0, $false, "28×55=", 1) | Out-Null
$t.Cell(4, 4).Range.Find.Execute("64×24=", $true, $false, $false, $false, $false, $true, 0, $false, "86×60=", 1) | Out-Null
$t.Cell(4, 5).Range.Find.Execute("54×33=", $true, $false, $false, $false, $false, $true, 0, $false, "88×32=", 1) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("80×56=", $true, $false, $false, $false, $false, $true, 0, $false, "40×43=", 1) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("22×36=", $true, $false, $false, $false, $false, $true, 0, $false, "24×36=", 1) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("88×37=", $true, $false, $false, $false, $false, $true, 0, $false, "33×82=", 1) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("85×57=", $true, $false, $false, $false, $false, $true, 0, $false, "97×13=", 1) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("45×84=", $true, $false, $false, $false, $false, $true, 0, $false, "28×25=", 1) | Out-Null
$t.Cell(6, 1).Range.Find.Execute("16×92=", $true, $false, $false, $false, $false, $true, 0, $false, "29×37=", 1) | Out-Null
$t.Cell(6, 2).Range.Find.Execute("69×41=", $true, $false, $false, $false, $false, $true, 0, $false, "22×23=", 1) | Out-Null
$t.Cell(6, 3).Range.Find.Execute("31×75=", $true, $false, $false, $false, $false, $true, 0, $false, "62×13=", 1) | Out-Null
$t.Cell(6, 4).Range.Find.Execute("51×77=", $true, $false, $false, $false, $false, $true, 0, $false, "59×67=", 1) | Out-Null
$t.Cell(6, 5).Range.Find.Execute("72×100=", $true, $false, $false, $false, $false, $true, 0, $false, "46×63=", 1) | Out-Null
$t.Cell(7, 1).Range.Find.Execute("24×98=", $true, $false, $false, $false, $false, $true, 0, $false, "79×60=", 1) | Out-Null
$t.Cell(7, 2).Range.Find.Execute("14×75=", $true, $false, $false, $false, $false, $true, 0, $false, "88×15=", 1) | Out-Null
$t.Cell(7, 3).Range.Find.Execute("50×46=", $true, $false, $false, $false, $false, $true, 0, $false, "68×40=", 1) | Out-Null
$t.Cell(7, 4).Range.Find.Execute("63×36=", $true, $false, $false, $false, $false, $true, 0, $false, "52×46=", 1) | Out-Null
$t.Cell(7, 5).Range.Find.Execute("64×50=", $true, $false, $false, $false, $false, $true, 0, $false, "48×30=", 1) | Out-Null
$t.Cell(8, 1).Range.Find.Execute("48×81=", $true, $false, $false, $false, $false, $true, 0, $false, "39×39=", 1) | Out-Null
$t.Cell(8, 2).Range.Find.Execute("83×55=", $true, $false, $false, $false, $false, $true, 0, $false, "60×91=", 1) | Out-Null
$t.Cell(8, 3).Range.Find.Execute("82×86=", $true, $false, $false, $false, $false, $true, 0, $false, "73×95=", 1) | Out-Null
$t.Cell(8, 4).Range.Find.Execute("15×25=", $true, $false, $false, $false, $false, $true, 0, $false, "35×65=", 1) | Out-Null
$t.Cell(8, 5).Range.Find.Execute("83×14=", $true, $false, $false, $false, $false, $true, 0, $false, "78×20=", 1) | Out-Null
$t.Cell(9, 1).Range.Find.Execute("57×52=", $true, $false, $false, $false, $false, $true, 0, $false, "94×49=", 1) | Out-Null
$t.Cell(9, 2).Range.Find.Execute("42×45=", $true, $false, $false, $false, $false, $true, 0, $false, "34×63=", 1) | Out-Null
$t.Cell(9, 3).Range.Find.Execute("37×54=", $true, $false, $false, $false, $false, $true, 0, $false, "11×73=", 1) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("84×17=", $true, $false, $false, $false, $false, $true, 0, $false, "71×59=", 1) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("76×52=", $true, $false, $false, $false, $false, $true, 0, $false, "57×23=", 1) | Out-Null
$t.Cell(10, 1).Range.Find.Execute("83×15=", $true, $false, $false, $false, $false, $true, 0, $false, "43×95=", 1) | Out-Null
$t.Cell(10, 2).Range.Find.Execute("72×35=", $true, $false, $false, $false, $false, $true, 0, $false, "26×80=", 1) | Out-Null
$t.Cell(10, 3).Range.Find.Execute("68×98=", $true, $false, $false, $false, $false, $true, 0, $false, "52×34=", 1) | Out-Null
$t.Cell(10, 4).Range.Find.Execute("46×40=", $true, $false, $false, $false, $false, $true, 0, $false, "57×46=", 1) | Out-Null
$t.Cell(10, 5).Range.Find.Execute("52×77=", $true, $false, $false, $false, $false, $true, 0, $false, "37×33=", 1) | Out-Null
$t.Cell(11, 1).Range.Find.Execute("87×59=", $true, $false, $false, $false, $false, $true, 0, $false, "76×67=", 1) | Out-Null
$t.Cell(11, 2).Range.Find.Execute("90×11=", $true, $false, $false, $false, $false, $true, 0, $false, "23×85=", 1) | Out-Null
$t.Cell(11, 3).Range.Find.Execute("25×40=", $true, $false, $false, $false, $false, $true, 0, $false, "47×12=", 1) | Out-Null
$t.Cell(11, 4).Range.Find.Execute("40×78=", $true, $false, $false, $false, $false, $true, 0, $false, "21×31=", 1) | Out-Null
$t.Cell(11, 5).Range.Find.Execute("32×35=", $true, $false, $false, $false, $false, $true, 0, $false, "12×61=", 1) | Out-Null
$t.Cell(12, 1).Range.Find.Execute("98×47=", $true, $false, $false, $false, $false, $true, 0, $false, "16×12=", 1) | Out-Null
$t.Cell(12, 2).Range.Find.Execute("69×53=", $true, $false, $false, $false, $false, $true, 0, $false, "79×46=", 1) | Out-Null
$t.Cell(12, 3).Range.Find.Execute("16×99=", $true, $false, $false, $false, $false, $true, 0, $false, "10×15=", 1) | Out-Null
$t.Cell(12, 4).Range.Find.Execute("84×54=", $true, $false, $false, $false, $false, $true, 0, $false, "16×21=", 1) | Out-Null
$t.Cell(12, 5).Range.Find.Execute("78×94=", $true, $false, $false, $false, $false, $true, 0, $false, "13×85=", 1) | Out-Null
$t.Cell(13, 1).Range.Find.Execute("76×55=", $true, $false, $false, $false, $false, $true, 0, $false, "44×96=", 1) | Out-Null
$t.Cell(13, 2).Range.Find.Execute("85×31=", $true, $false, $false, $false, $false, $true, 0, $false, "34×73=", 1) | Out-Null
$t.Cell(13, 3).Range.Find.Execute("20×34=", $true, $false, $false, $false, $false, $true, 0, $false, "57×62=", 1) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("13×37=", $true, $false, $false, $false, $false, $true, 0, $false, "20×92=", 1) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("86×43=", $true, $false, $false, $false, $false, $true, 0, $false, "34×63=", 1) | Out-Null
$t.Cell(14, 1).Range.Find.Execute("30×53=", $true, $false, $false, $false, $false, $true, 0, $false, "83×98=", 1) | Out-Null
$t.Cell(14, 2).Range.Find.Execute("34×20=", $true, $false, $false, $false, $false, $true, 0, $false, "30×86=", 1) | Out-Null
$t.Cell(14, 3).Range.Find.Execute("54×95=", $true, $false, $false, $false, $false, $true, 0, $false, "92×12=", 1) | Out-Null
$t.Cell(14, 4).Range.Find.Execute("99×24=", $true, $false, $false, $false, $false, $true, 0, $false, "41×79=", 1) | Out-Null
$t.Cell(14, 5).Range.Find.Execute("82×52=", $true, $false, $false, $false, $false, $true, 0, $false, "67×97=", 1) | Out-Null
$t.Cell(15, 1).Range.Find.Execute("14×41=", $true, $false, $false, $false, $false, $true, 0, $false, "100×60=", 1) | Out-Null
$t.Cell(15, 2).Range.Find.Execute("87×52=", $true, $false, $false, $false, $false, $true, 0, $false, "62×69=", 1) | Out-Null
$t.Cell(15, 3).Range.Find.Execute("92×16=", $true, $false, $false, $false, $false, $true, 0, $false, "82×99=", 1) | Out-Null
$t.Cell(15, 4).Range.Find.Execute("42×16=", $true, $false, $false, $false, $false, $true, 0, $false, "59×51=", 1) | Out-Null
$t.Cell(15, 5).Range.Find.Execute("34×51=", $true, $false, $false, $false, $false, $true, 0, $false, "82×75=", 1) | Out-Null
$t.Cell(16, 1).Range.Find.Execute("16×30=", $true, $false, $false, $false, $false, $true, 0, $false, "51×35=", 1) | Out-Null
$t.Cell(16, 2).Range.Find.Execute("15×44=", $true, $false, $false, $false, $false, $true, 0, $false, "66×75=", 1) | Out-Null
$t.Cell(16, 3).Range.Find.Execute("34×25=", $true, $false, $false, $false, $false, $true, 0, $false, "24×81=", 1) | Out-Null
$t.Cell(16, 4).Range.Find.Execute("20×24=", $true, $false, $false, $false, $false, $true, 0, $false, "100×93=", 1) | Out-Null
$t.Cell(16, 5).Range.Find.Execute("61×83=", $true, $false, $false, $false, $false, $true, 0, $false, "58×84=", 1) | Out-Null
$t.Cell(17, 1).Range.Find.Execute("43×25=", $true, $false, $false, $false, $false, $true, 0, $false, "30×56=", 1) | Out-Null
$t.Cell(17, 2).Range.Find.Execute("58×39=", $true, $false, $false, $false, $false, $true, 0, $false, "48×47=", 1) | Out-Null
$t.Cell(17, 3).Range.Find.Execute("16×34=", $true, $false, $false, $false, $false, $true, 0, $false, "15×11=", 1) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("11×56=", $true, $false, $false, $false, $false, $true, 0, $false, "80×42=", 1) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("10×55=", $true, $false, $false, $false, $false, $true, 0, $false, "27×25=", 1) | Out-Null
$t.Cell(18, 1).Range.Find.Execute("97×83=", $true, $false, $false, $false, $false, $true, 0, $false, "46×14=", 1) | Out-Null
$t.Cell(18, 2).Range.Find.Execute("84×84=", $true, $false, $false, $false, $false, $true, 0, $false, "56×82=", 1) | Out-Null
$t.Cell(18, 3).Range.Find.Execute("76×14=", $true, $false, $false, $false, $false, $true, 0, $false, "93×11=", 1) | Out-Null
$t.Cell(18, 4).Range.Find.Execute("83×73=", $true, $false, $false, $false, $false, $true, 0, $false, "92×21=", 1) | Out-Null
$t.Cell(18, 5).Range.Find.Execute("26×74=", $true, $false, $false, $false, $false, $true, 0, $false, "63×86=", 1) | Out-Null
$t.Cell(19, 1).Range.Find.Execute("38×98=", $true, $false, $false, $false, $false, $true, 0, $false, "62×35=", 1) | Out-Null
$t.Cell(19, 2).Range.Find.Execute("17×13=", $true, $false, $false, $false, $false, $true, 0, $false, "66×100=", 1) | Out-Null
$t.Cell(19, 3).Range.Find.Execute("44×77=", $true, $false, $false, $false, $false, $true, 0, $false, "56×79=", 1) | Out-Null
$t.Cell(19, 4).Range.Find.Execute("51×43=", $true, $false, $false, $false, $false, $true, 0, $false, "93×90=", 1) | Out-Null
$t.Cell(19, 5).Range.Find.Execute("23×75=", $true, $false, $false, $false, $false, $true, 0, $false, "14×78=", 1) | Out-Null
$t.Cell(20, 1).Range.Find.Execute("28×48=", $true, $false, $false, $false, $false, $true, 0, $false, "87×41=", 1) | Out-Null
$t.Cell(20, 2).Range.Find.Execute("16×96=", $true, $false, $false, $false, $false, $true, 0, $false, "93×37=", 1) | Out-Null
$t.Cell(20, 3).Range.Find.Execute("32×94=", $true, $false, $false, $false, $false, $true, 0, $false, "75×37=", 1) | Out-Null
$t.Cell(20, 4).Range.Find.Execute("70×19=", $true, $false, $false, $false, $false, $true, 0, $false, "96×11=", 1) | Out-Null
$t.Cell(20, 5).Range.Find.Execute("53×49=", $true, $false, $false, $false, $false, $true, 0, $false, "21×26=", 1) | Out-Null
